$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G3").Value = "2016-08-30 00:46:42"

$wsZhCn.Range("H3").Value = "2016-08-30 00:46:37"
$wsZhCn.Range("K3").Value = "2016-08-30 00:46:56"

$wsDeDe.Range("H3").Value = "2016-08-30 00:46:42"
$wsDeDe.Range("K3").Value = "2016-08-30 00:47:09"
